$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.08897633333333334
$ws.Range("H2").Value = 0.266929
$ws.Range("I2").Value = 0.01166961909325238
$ws.Range("J2").Value = 0.01190636343604672
$ws.Range("M2").Value = 0.09788266666666667
$ws.Range("N2").Value = 0.293648
$ws.Range("Q2").Value = 0.00870924077688889
$ws.Range("R2").Value = 0.07838316699200001
$ws.Range("S2").Value = 0.01166961909325238
$ws.Range("T2").Value = 0.01190636343604672

$ws.Range("I3").Value = 0.4763688563004725
$ws.Range("J3").Value = 0.4860330647816007
$ws.Range("M3").Value = 0.09788266666666667
$ws.Range("N3").Value = 0.293648
$ws.Range("Q3").Value = 0.3555224069422223
$ws.Range("R3").Value = 3.19970166248
$ws.Range("S3").Value = 0.4763688563004725
$ws.Range("T3").Value = 0.4860330647816007

$ws.Range("G4").Value = 1.962678
$ws.Range("H4").Value = 5.888033999999999
$ws.Range("I4").Value = 0.2574134469769833
$ws.Range("J4").Value = 0.262635654903738
$ws.Range("M4").Value = 0.09788266666666667
$ws.Range("N4").Value = 0.293648
$ws.Range("Q4").Value = 0.192112156448
$ws.Range("R4").Value = 1.729009408032
$ws.Range("S4").Value = 0.2574134469769833
$ws.Range("T4").Value = 0.262635654903738

$ws.Range("G5").Value = 0.45482
$ws.Range("H5").Value = 0.90964
$ws.Range("I5").Value = 0.0596515495430588
$ws.Range("J5").Value = 0.04057447649362016
$ws.Range("M5").Value = 0.09788266666666667
$ws.Range("N5").Value = 0.293648
$ws.Range("Q5").Value = 0.04451899445333334
$ws.Range("R5").Value = 0.26711396672
$ws.Range("S5").Value = 0.0596515495430588
$ws.Range("T5").Value = 0.04057447649362016

$ws.Range("G6").Value = 1.486010666666667
$ws.Range("H6").Value = 4.458032
$ws.Range("I6").Value = 0.194896528086233
$ws.Range("J6").Value = 0.1988504403849946
$ws.Range("M6").Value = 0.09788266666666667
$ws.Range("N6").Value = 0.293648
$ws.Range("Q6").Value = 0.1454546867484445
$ws.Range("R6").Value = 1.309092180736
$ws.Range("S6").Value = 0.194896528086233
$ws.Range("T6").Value = 0.1988504403849946
